# Update odds values in Sheet1 to match the latest FlashScore data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 2.38
$ws.Range("L2").Value = 5.5
$ws.Range("N2").Value = 8.5
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.75
$ws.Range("W2").Value = 6
$ws.Range("AC2").Value = 8.5
$ws.Range("AG2").Value = 401
$ws.Range("AI2").Value = 26
$ws.Range("AO2").Value = 9
$ws.Range("AQ2").Value = 29
$ws.Range("AX2").Value = 29

# Row 3
$ws.Range("G3").Value = 1.95
$ws.Range("I3").Value = 4.1
$ws.Range("W3").Value = 5.5
$ws.Range("AM3").Value = 51
$ws.Range("AS3").Value = 251
$ws.Range("AX3").Value = 26

# Row 6
$ws.Range("P6").Value = 4
